# Update weekly Fruta/Hortaliza price data for Femacal de La Calera - Breva.
# The underlying data rows (2,4,5,6,7,8,9) are being reordered: each row's
# Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) values are shuffled
# among the rows (row 3 is untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that move, keyed by row.
$rows = @(2, 4, 5, 6, 7, 8, 9)
$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{
        D = $ws.Range("D$r").Value2
        M = $ws.Range("M$r").Value2
        N = $ws.Range("N$r").Value2
        O = $ws.Range("O$r").Value2
        P = $ws.Range("P$r").Value2
        S = $ws.Range("S$r").Value2
    }
}

# Destination row -> source row (i.e. destination row now holds the values
# that used to live in the source row).
$mapping = @{
    2 = 7
    4 = 6
    5 = 4
    6 = 2
    7 = 5
    8 = 9
    9 = 8
}

foreach ($dest in $mapping.Keys) {
    $src = $mapping[$dest]
    $vals = $before[$src]
    $ws.Range("D$dest").Value = $vals.D
    $ws.Range("M$dest").Value = $vals.M
    $ws.Range("N$dest").Value = $vals.N
    $ws.Range("O$dest").Value = $vals.O
    $ws.Range("P$dest").Value = $vals.P
    $ws.Range("S$dest").Value = $vals.S
}
